# Scheduled runner update: refresh currentAveragePrice / Leve profit figures
# across several job sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR) to match the
# latest market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1295.4546
$ws.Range("I18").Value = 1225
$ws.Range("K18").Value = 1225
$ws.Range("M18").Value = -941

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2212.8125
$ws.Range("I41").Value = 2350.5
$ws.Range("J41").Value = 1983.3334
$ws.Range("K41").Value = 2350.5
$ws.Range("L41").Value = 1983.3334
$ws.Range("M41").Value = -1910.5
$ws.Range("N41").Value = -2863.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3732.08
$ws.Range("I64").Value = 3805.8823
$ws.Range("J64").Value = 3575.25
$ws.Range("K64").Value = 3805.8823
$ws.Range("L64").Value = 3575.25
$ws.Range("M64").Value = -3557.8823
$ws.Range("N64").Value = -4071.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3732.08
$ws.Range("I67").Value = 3805.8823
$ws.Range("J67").Value = 3575.25
$ws.Range("K67").Value = 3805.8823
$ws.Range("L67").Value = 3575.25
$ws.Range("M67").Value = -2947.8823
$ws.Range("N67").Value = -5291.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 576.17145
$ws.Range("J129").Value = 878.73334
$ws.Range("L129").Value = 2636.20002
$ws.Range("N129").Value = -12636.20002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3791346.5
$ws.Range("I132").Value = 4697767
$ws.Range("J132").Value = 5708.8237
$ws.Range("K132").Value = 14093301
$ws.Range("L132").Value = 17126.4711
$ws.Range("M132").Value = -14090771
$ws.Range("N132").Value = -22186.4711

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1196.8572
$ws.Range("I137").Value = 882.74286
$ws.Range("J137").Value = 1589.5
$ws.Range("K137").Value = 2648.22858
$ws.Range("L137").Value = 4768.5
$ws.Range("M137").Value = -98.22857999999997
$ws.Range("N137").Value = -9868.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 5006
$ws.Range("I16").Value = 5006
$ws.Range("K16").Value = 5006
$ws.Range("M16").Value = -4719

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5132.561
$ws.Range("I32").Value = 4382.485
$ws.Range("J32").Value = 8226.625
$ws.Range("K32").Value = 4382.485
$ws.Range("L32").Value = 8226.625
$ws.Range("M32").Value = -4095.485
$ws.Range("N32").Value = -8800.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 249800
$ws.Range("J48").Value = 249800
$ws.Range("L48").Value = 249800
$ws.Range("N48").Value = -250568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1070.9615
$ws.Range("I74").Value = 749.1739
$ws.Range("J74").Value = 3538
$ws.Range("K74").Value = 749.1739
$ws.Range("L74").Value = 3538
$ws.Range("M74").Value = 124.8261
$ws.Range("N74").Value = -5286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1070.9615
$ws.Range("I77").Value = 749.1739
$ws.Range("J77").Value = 3538
$ws.Range("K77").Value = 3745.8695
$ws.Range("L77").Value = 17690
$ws.Range("M77").Value = 622.1305000000002
$ws.Range("N77").Value = -26426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2309.8965
$ws.Range("I132").Value = 2349
$ws.Range("J132").Value = 2268
$ws.Range("K132").Value = 7047
$ws.Range("L132").Value = 6804
$ws.Range("M132").Value = -4517
$ws.Range("N132").Value = -11864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 36147.8
$ws.Range("J117").Value = 36147.8
$ws.Range("L117").Value = 36147.8
$ws.Range("N117").Value = -45325.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3275.6
$ws.Range("I134").Value = 1140.3784
$ws.Range("J134").Value = 9352.77
$ws.Range("K134").Value = 3421.1352
$ws.Range("L134").Value = 28058.31
$ws.Range("M134").Value = -886.1352000000002
$ws.Range("N134").Value = -33128.31

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2728.8235
$ws.Range("I31").Value = 2899.4285
$ws.Range("J31").Value = 1932.6666
$ws.Range("K31").Value = 2899.4285
$ws.Range("L31").Value = 1932.6666
$ws.Range("M31").Value = -2604.4285
$ws.Range("N31").Value = -2522.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2728.8235
$ws.Range("I34").Value = 2899.4285
$ws.Range("J34").Value = 1932.6666
$ws.Range("K34").Value = 2899.4285
$ws.Range("L34").Value = 1932.6666
$ws.Range("M34").Value = -2697.4285
$ws.Range("N34").Value = -2336.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 883.61816
$ws.Range("I58").Value = 847.73334
$ws.Range("K58").Value = 847.73334
$ws.Range("M58").Value = -644.73334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3061.7458
$ws.Range("I132").Value = 3696.825
$ws.Range("J132").Value = 1724.7368
$ws.Range("K132").Value = 11090.475
$ws.Range("L132").Value = 5174.2104
$ws.Range("M132").Value = -8560.474999999999
$ws.Range("N132").Value = -10234.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9805114
$ws.Range("I134").Value = 1216.6945
$ws.Range("J134").Value = 33334468
$ws.Range("K134").Value = 3650.0835
$ws.Range("L134").Value = 100003404
$ws.Range("M134").Value = -1115.0835
$ws.Range("N134").Value = -100008474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 883.61816
$ws.Range("I136").Value = 847.73334
$ws.Range("K136").Value = 2543.20002
$ws.Range("M136").Value = 6.799979999999778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2463.8667
$ws.Range("I136").Value = 1870
$ws.Range("J136").Value = 3142.5715
$ws.Range("K136").Value = 5610
$ws.Range("L136").Value = 9427.7145
$ws.Range("M136").Value = -510
$ws.Range("N136").Value = -19627.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 12000
$ws.Range("J117").Value = 12000
$ws.Range("L117").Value = 12000
$ws.Range("N117").Value = -18884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6758465.5
$ws.Range("J122").Value = 27779616
$ws.Range("L122").Value = 83338848
$ws.Range("N122").Value = -83343748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3823.4092
$ws.Range("I132").Value = 5071.8887
$ws.Range("J132").Value = 2959.077
$ws.Range("K132").Value = 15215.6661
$ws.Range("L132").Value = 8877.231
$ws.Range("M132").Value = -12685.6661
$ws.Range("N132").Value = -13937.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 33293.168
$ws.Range("J135").Value = 29951.8
$ws.Range("L135").Value = 29951.8
$ws.Range("N135").Value = -40091.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 24691.555
$ws.Range("J136").Value = 24691.555
$ws.Range("L136").Value = 74074.66500000001
$ws.Range("N136").Value = -79174.66500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 50000
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("N16").Value = -50584

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1945.3286
$ws.Range("I132").Value = 1792.2
$ws.Range("K132").Value = 5376.6
$ws.Range("M132").Value = -2846.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 538.8889
$ws.Range("I136").Value = 470.58823
$ws.Range("K136").Value = 1411.76469
$ws.Range("M136").Value = 1138.23531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
